$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("lte")

$ws.Range("C2").Value = "L2100cell"
$ws.Range("B3").Value = "ERBS_555"
$ws.Range("C3").Value = "L1800"
$ws.Range("B4").Value = "ERBS_777"
$ws.Range("C4").Value = "cell3"
